# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45181 (2023-09-12) to 45182 (2023-09-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data rows are 2 through 482 (column C holds the "Förändrad" date)
$lastRow = 482

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45182
